# ------------------------------------------------------------------
# Applies the "Update OnSSET input files" edit:
#  - ScenarioInfo: splits SAPV_capital_cost into 5 columns
#      (SAPV_capital_cost_Inf/_1kW/_100W/_50W/_20W), shifting the
#      following Diesel columns 4 to the right.
#  - ScenarioParameters: splits SAPVCapitalCost into 5 columns
#      (SAPVCapitalCostInf/1kW/100W/50W/20W) with new cost figures,
#      updates MGHydroCapitalCost / MGWindCapitalCost / MGDieselCapitalCost
#      / SADieselCapitalCost values.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet: ScenarioInfo  (header columns shift: insert 4 new cols after AP)
# =========================================================================
$wsInfo = $wb.Worksheets.Item("ScenarioInfo")

# Insert 4 new (blank) columns where AQ used to be - this pushes the old
# AQ:BI block to AU:BM automatically.
$wsInfo.Columns("AQ:AT").Insert()

# Rename the (now split) capital-cost header and label the 4 new columns.
$wsInfo.Range("AP1").Value = "SAPV_capital_cost_Inf"
$wsInfo.Range("AQ1").Value = "SAPV_capital_cost_1kW"
$wsInfo.Range("AR1").Value = "SAPV_capital_cost_100W"
$wsInfo.Range("AS1").Value = "SAPV_capital_cost_50W"
$wsInfo.Range("AT1").Value = "SAPV_capital_cost_20W"

# Every data row in this sheet repeats a single constant across the whole
# row (row 2 -> 0, row 3 -> 1, row 4 -> 2, row 5 -> 3, row 6 -> 4), so the
# 4 freshly inserted (blank) cells just need that same row constant.
$wsInfo.Range("AQ2:AT2").Value = 0
$wsInfo.Range("AQ3:AT3").Value = 1
$wsInfo.Range("AQ4:AT4").Value = 2
$wsInfo.Range("AQ5:AT5").Value = 3
$wsInfo.Range("AQ6:AT6").Value = 4

# =========================================================================
# Sheet: ScenarioParameters  (header columns shift: insert 4 new cols
# after AS, plus several value updates)
# =========================================================================
$wsParams = $wb.Worksheets.Item("ScenarioParameters")

# Insert 4 new (blank) columns where AT used to be - pushes old AT:BL to
# AX:BP automatically.
$wsParams.Columns("AT:AW").Insert()

# Rename the (now split) capital-cost header and label the 4 new columns.
$wsParams.Range("AS1").Value = "SAPVCapitalCostInf"
$wsParams.Range("AT1").Value = "SAPVCapitalCost1kW"
$wsParams.Range("AU1").Value = "SAPVCapitalCost100W"
$wsParams.Range("AV1").Value = "SAPVCapitalCost50W"
$wsParams.Range("AW1").Value = "SAPVCapitalCost20W"

# --- MGHydroCapitalCost (column X) updated to 3000 for every row ---
$wsParams.Range("X2:X6").Value = 3000

# --- MGWindCapitalCost (column AF) updated per row ---
$wsParams.Range("AF2").Value = 2800
$wsParams.Range("AF3").Value = 2213.08840413318
$wsParams.Range("AF4").Value = 2142.365097588978
$wsParams.Range("AF5").Value = 2142.365097588978
$wsParams.Range("AF6").Value = 2142.365097588978

# --- Split SAPV capital-cost figures (AS..AW) ---
$wsParams.Range("AS2").Value = 6950
$wsParams.Range("AT2").Value = 4470
$wsParams.Range("AU2").Value = 6380
$wsParams.Range("AV2").Value = 8780
$wsParams.Range("AW2").Value = 9620

$wsParams.Range("AS3").Value = 5340.618347203406
$wsParams.Range("AT3").Value = 3434.901296690536
$wsParams.Range("AU3").Value = 4902.61079930327
$wsParams.Range("AV3").Value = 6746.853106251209
$wsParams.Range("AW3").Value = 7392.337913682988

$wsParams.Range("AS4").Value = 4862.444358428489
$wsParams.Range("AT4").Value = 3127.356299593574
$wsParams.Range("AU4").Value = 4463.653957809173
$wsParams.Range("AV4").Value = 6142.771434101026
$wsParams.Range("AW4").Value = 6730.462550803174

$wsParams.Range("AS5").Value = 4862.444358428489
$wsParams.Range("AT5").Value = 3127.356299593574
$wsParams.Range("AU5").Value = 4463.653957809173
$wsParams.Range("AV5").Value = 6142.771434101026
$wsParams.Range("AW5").Value = 6730.462550803174

$wsParams.Range("AS6").Value = 4862.444358428489
$wsParams.Range("AT6").Value = 3127.356299593574
$wsParams.Range("AU6").Value = 4463.653957809173
$wsParams.Range("AV6").Value = 6142.771434101026
$wsParams.Range("AW6").Value = 6730.462550803174

# --- MGDieselCapitalCost (now column BB, was AX) -> 261 for every row ---
$wsParams.Range("BB2:BB6").Value = 261

# --- SADieselCapitalCost (now column BI, was BE) -> 261 for every row ---
$wsParams.Range("BI2:BI6").Value = 261

# =========================================================================
# Workbook view: keep "SpecsData" (3rd tab) as the active sheet, matching
# the saved activeTab index.
# =========================================================================
$wb.Worksheets.Item("SpecsData").Activate()

# =========================================================================
# Cosmetic sheet-view / page-setup tweaks for ScenarioInfo that are
# reachable through the COM surface.
# =========================================================================
$wsInfo.Activate()
$winInfo = $excel.ActiveWindow
$winInfo.Zoom = 100
$wsInfo.Range("BJ1").Select()

$wsInfo.PageSetup.Orientation = 1
$wsInfo.PageSetup.PaperSize = 9
$wsInfo.PageSetup.HeaderMargin = 36.850393700787386
$wsInfo.PageSetup.FooterMargin = 36.850393700787386

# Re-activate SpecsData so the workbook's active tab matches the target.
$wb.Worksheets.Item("SpecsData").Activate()
